# Scheduled runner update: refresh market-price driven columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Sheets workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 29
$ws.Range("H29").Value = 842.6667
$ws.Range("I29").Value = 32.4
$ws.Range("J29").Value = 4894
$ws.Range("K29").Value = 97.19999999999999
$ws.Range("L29").Value = 14682
$ws.Range("M29").Value = 183.8
$ws.Range("N29").Value = -15244

# ALC row 86
$ws.Range("H86").Value = 2410
$ws.Range("I86").Value = 2325
$ws.Range("J86").Value = 2750
$ws.Range("K86").Value = 2325
$ws.Range("L86").Value = 2750
$ws.Range("M86").Value = -1202
$ws.Range("N86").Value = -4996

# ALC row 89
$ws.Range("H89").Value = 2410
$ws.Range("I89").Value = 2325
$ws.Range("J89").Value = 2750
$ws.Range("K89").Value = 11625
$ws.Range("L89").Value = 13750
$ws.Range("M89").Value = -6009
$ws.Range("N89").Value = -24982

# ALC row 92
$ws.Range("H92").Value = 14933.25
$ws.Range("I92").Value = 16652.285
$ws.Range("J92").Value = 2900
$ws.Range("K92").Value = 16652.285
$ws.Range("L92").Value = 2900
$ws.Range("M92").Value = -15404.285
$ws.Range("N92").Value = -5396

# ALC row 138
$ws.Range("H138").Value = 2375.5715
$ws.Range("I138").Value = 1398.9796
$ws.Range("J138").Value = 4084.6072
$ws.Range("K138").Value = 4196.9388
$ws.Range("L138").Value = 12253.8216
$ws.Range("M138").Value = 943.0612000000001
$ws.Range("N138").Value = -22533.8216

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 17116.7
$ws.Range("I32").Value = 20028
$ws.Range("J32").Value = 8382.799999999999
$ws.Range("K32").Value = 20028
$ws.Range("L32").Value = 8382.799999999999
$ws.Range("M32").Value = -19741
$ws.Range("N32").Value = -8956.799999999999

# ARM row 61
$ws.Range("H61").Value = 4311.2104
$ws.Range("I61").Value = 3409.4167
$ws.Range("J61").Value = 5857.143
$ws.Range("K61").Value = 3409.4167
$ws.Range("L61").Value = 5857.143
$ws.Range("M61").Value = -3197.4167
$ws.Range("N61").Value = -6281.143

# ARM row 88
$ws.Range("H88").Value = 2621.2
$ws.Range("I88").Value = 2553
$ws.Range("J88").Value = 2666.6667
$ws.Range("K88").Value = 2553
$ws.Range("L88").Value = 2666.6667
$ws.Range("M88").Value = -2147
$ws.Range("N88").Value = -3478.6667

# ARM row 91
$ws.Range("H91").Value = 2621.2
$ws.Range("I91").Value = 2553
$ws.Range("J91").Value = 2666.6667
$ws.Range("K91").Value = 2553
$ws.Range("L91").Value = 2666.6667
$ws.Range("M91").Value = -1149
$ws.Range("N91").Value = -5474.6667

# ARM row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# ARM row 131
$ws.Range("H131").Value = 54678.5
$ws.Range("J131").Value = 54678.5
$ws.Range("L131").Value = 54678.5
$ws.Range("N131").Value = -64758.5

# ARM row 136
$ws.Range("H136").Value = 4311.2104
$ws.Range("I136").Value = 3409.4167
$ws.Range("J136").Value = 5857.143
$ws.Range("K136").Value = 10228.2501
$ws.Range("L136").Value = 17571.429
$ws.Range("M136").Value = -7678.250100000001
$ws.Range("N136").Value = -22671.429

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86
$ws.Range("H86").Value = 289859
$ws.Range("I86").Value = 5801.2
$ws.Range("J86").Value = 1000003.5
$ws.Range("K86").Value = 5801.2
$ws.Range("L86").Value = 1000003.5
$ws.Range("M86").Value = -4678.2
$ws.Range("N86").Value = -1002249.5

# BSM row 89
$ws.Range("H89").Value = 289859
$ws.Range("I89").Value = 5801.2
$ws.Range("J89").Value = 1000003.5
$ws.Range("K89").Value = 29006
$ws.Range("L89").Value = 5000017.5
$ws.Range("M89").Value = -23390
$ws.Range("N89").Value = -5011249.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58
$ws.Range("H58").Value = 1030281.7
$ws.Range("I58").Value = 1235777.1
$ws.Range("J58").Value = 2804.5
$ws.Range("K58").Value = 1235777.1
$ws.Range("L58").Value = 2804.5
$ws.Range("M58").Value = -1235574.1
$ws.Range("N58").Value = -3210.5

# CRP row 136
$ws.Range("H136").Value = 1030281.7
$ws.Range("I136").Value = 1235777.1
$ws.Range("J136").Value = 2804.5
$ws.Range("K136").Value = 3707331.3
$ws.Range("L136").Value = 8413.5
$ws.Range("M136").Value = -3704781.3
$ws.Range("N136").Value = -13513.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 40
$ws.Range("H40").Value = 290
$ws.Range("I40").Value = 186
$ws.Range("J40").Value = 550
$ws.Range("K40").Value = 744
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -675
$ws.Range("N40").Value = -2338

# CUL row 68
$ws.Range("H68").Value = 797.125
$ws.Range("I68").Value = 1020
$ws.Range("J68").Value = 722.8333
$ws.Range("K68").Value = 3060
$ws.Range("L68").Value = 2168.4999
$ws.Range("M68").Value = -2249
$ws.Range("N68").Value = -3790.4999

# CUL row 69
$ws.Range("H69").Value = 1591.6666
$ws.Range("J69").Value = 1810
$ws.Range("L69").Value = 5430
$ws.Range("N69").Value = -7052

# CUL row 71
$ws.Range("H71").Value = 797.125
$ws.Range("I71").Value = 1020
$ws.Range("J71").Value = 722.8333
$ws.Range("K71").Value = 9180
$ws.Range("L71").Value = 6505.4997
$ws.Range("M71").Value = -5124
$ws.Range("N71").Value = -14617.4997

# CUL row 72
$ws.Range("H72").Value = 1591.6666
$ws.Range("J72").Value = 1810
$ws.Range("L72").Value = 16290
$ws.Range("N72").Value = -24402

# CUL row 80
$ws.Range("H80").Value = 5419.615
$ws.Range("I80").Value = 6985.7144
$ws.Range("J80").Value = 3592.5
$ws.Range("K80").Value = 20957.1432
$ws.Range("L80").Value = 10777.5
$ws.Range("M80").Value = -20021.1432
$ws.Range("N80").Value = -12649.5

# CUL row 83
$ws.Range("H83").Value = 5419.615
$ws.Range("I83").Value = 6985.7144
$ws.Range("J83").Value = 3592.5
$ws.Range("K83").Value = 62871.4296
$ws.Range("L83").Value = 32332.5
$ws.Range("M83").Value = -58191.4296
$ws.Range("N83").Value = -41692.5

# CUL row 109
$ws.Range("H109").Value = 3392
$ws.Range("I109").Value = 566
$ws.Range("J109").Value = 6218
$ws.Range("K109").Value = 1698
$ws.Range("L109").Value = 18654
$ws.Range("M109").Value = -658
$ws.Range("N109").Value = -20734

# CUL row 122
$ws.Range("H122").Value = 723.7646999999999
$ws.Range("J122").Value = 975.375
$ws.Range("L122").Value = 8778.375
$ws.Range("N122").Value = -13678.375

# CUL row 131
$ws.Range("H131").Value = 980.95
$ws.Range("I131").Value = 583.3333
$ws.Range("J131").Value = 993.24744
$ws.Range("K131").Value = 1749.9999
$ws.Range("L131").Value = 2979.74232
$ws.Range("M131").Value = 3290.0001
$ws.Range("N131").Value = -13059.74232

$ws = $wb.Worksheets.Item("GSM")
# GSM row 51
$ws.Range("H51").Value = 15092.2
$ws.Range("J51").Value = 15092.2
$ws.Range("L51").Value = 15092.2
$ws.Range("N51").Value = -16110.2

# GSM row 70
$ws.Range("H70").Value = 6564.409
$ws.Range("I70").Value = 5831.385
$ws.Range("K70").Value = 5831.385
$ws.Range("M70").Value = -5561.385

# GSM row 73
$ws.Range("H73").Value = 6564.409
$ws.Range("I73").Value = 5831.385
$ws.Range("K73").Value = 5831.385
$ws.Range("M73").Value = -4895.385

# GSM row 131
$ws.Range("H131").Value = 30664
$ws.Range("J131").Value = 30664
$ws.Range("L131").Value = 30664
$ws.Range("N131").Value = -40744

# GSM row 139
$ws.Range("H139").Value = 152599.33
$ws.Range("J139").Value = 152599.33
$ws.Range("L139").Value = 152599.33
$ws.Range("N139").Value = -162879.33

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 3879.7778
$ws.Range("I7").Value = 4521.6
$ws.Range("J7").Value = 3077.5
$ws.Range("K7").Value = 4521.6
$ws.Range("L7").Value = 3077.5
$ws.Range("M7").Value = -4409.6
$ws.Range("N7").Value = -3301.5

# LTW row 126
$ws.Range("H126").Value = 3879.7778
$ws.Range("I126").Value = 4521.6
$ws.Range("J126").Value = 3077.5
$ws.Range("K126").Value = 13564.8
$ws.Range("L126").Value = 9232.5
$ws.Range("M126").Value = -11094.8
$ws.Range("N126").Value = -14172.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 125
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

Write-Host "Applied scheduled market-data refresh."
